# Bmp6-Bmpr1b.xlsx: refresh the LR-pair edge table with new TPM-based
# values and extend it from 6 rows to the full 3x3 sender/target grid
# (ECs / FAPs / MuSCs), matching the updated natmi TPM recompute.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs", "Bmp6", "Bmpr1b", "ECs", 3, 1, 7.739652666666667, 23.218958, 0.6488398532974882, 0.6488398532974882, 1, 0.3333333333333333, 0.09074700000000001, 0.272241, 0.04501122713837679, 0.04501122713837678, 0.7023502605420001, 6.321152344878, 0.02920507801320432, 0.02920507801320431),
    @("ECs", "Bmp6", "Bmpr1b", "FAPs", 3, 1, 7.739652666666667, 23.218958, 0.6488398532974882, 0.6488398532974882, 3, 1, 1.394935666666666, 4.184806999999999, 0.6918990835593063, 0.6918990835593063, 10.79631755234511, 97.16685797110598, 0.4489316998732868, 0.4489316998732868),
    @("ECs", "Bmp6", "Bmpr1b", "MuSCs", 3, 1, 7.739652666666667, 23.218958, 0.6488398532974882, 0.6488398532974882, 3, 1, 0.5304143333333333, 1.591243, 0.263089689302317, 0.263089689302317, 4.105222709421556, 36.947004384794, 0.1707030754109971, 0.1707030754109971),
    @("FAPs", "Bmp6", "Bmpr1b", "ECs", 3, 1, 3.641794, 10.925382, 0.3053032463428815, 0.3053032463428815, 1, 0.3333333333333333, 0.09074700000000001, 0.272241, 0.04501122713837679, 0.04501122713837678, 0.3304818801180001, 2.974336921062001, 0.01374207376722324, 0.01374207376722324),
    @("FAPs", "Bmp6", "Bmpr1b", "FAPs", 3, 1, 3.641794, 10.925382, 0.3053032463428815, 0.3053032463428815, 3, 1, 1.394935666666666, 4.184806999999999, 0.6918990835593063, 0.6918990835593063, 5.080068341252666, 45.720615071274, 0.2112390363523208, 0.2112390363523208),
    @("FAPs", "Bmp6", "Bmpr1b", "MuSCs", 3, 1, 3.641794, 10.925382, 0.3053032463428815, 0.3053032463428815, 3, 1, 0.5304143333333333, 1.591243, 0.263089689302317, 0.263089689302317, 1.931659736647334, 17.384937629826, 0.08032213622333745, 0.08032213622333745),
    @("MuSCs", "Bmp6", "Bmpr1b", "ECs", 3, 1, 0.5470016666666667, 1.641005, 0.04585690035963046, 0.04585690035963046, 1, 0.3333333333333333, 0.09074700000000001, 0.272241, 0.04501122713837679, 0.04501122713837678, 0.04963876024500001, 0.4467488422050001, 0.002064075357949239, 0.002064075357949239),
    @("MuSCs", "Bmp6", "Bmpr1b", "FAPs", 3, 1, 0.5470016666666667, 1.641005, 0.04585690035963046, 0.04585690035963046, 3, 1, 1.394935666666666, 4.184806999999999, 0.6918990835593063, 0.6918990835593063, 0.7630321345594443, 6.867289211034999, 0.03172834733369874, 0.03172834733369874),
    @("MuSCs", "Bmp6", "Bmpr1b", "MuSCs", 3, 1, 0.5470016666666667, 1.641005, 0.04585690035963046, 0.04585690035963046, 3, 1, 0.5304143333333333, 1.591243, 0.263089689302317, 0.263089689302317, 0.2901375243572222, 2.611237719215, 0.01206447766798249, 0.01206447766798249)
)
$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $data[$i]
    $r = $startRow + $i
    for ($c = 1; $c -le $row.Count; $c++) {
        $ws.Cells.Item($r, $c).Value = $row[$c - 1]
    }
}
